$d = $word.ActiveDocument

# 1. Merge the date/revision line runs into a single run's text.
$d.Content.Find.Execute(
    "11/0",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "11/08",
    2)

$d.Content.Find.Execute(
    "8/2020 Revised by Kim Nguyen",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "/08/2020 Revised by Kim Nguyen",
    2)

# 2. Merge "Center" + " for Information Assurance..." into a single run
#    and drop the spell-check markers around "Center".
$d.Content.Find.Execute(
    "Center for Information Assurance (CIAE) @City University of Seattle (",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "Center for Information Assurance (CIAE) @City University of Seattle (",
    2)

# 3. Merge "HOP0" + "6" + " assignment (..." into a single run.
$d.Content.Find.Execute(
    "HOP06 assignment",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "HOP06 assignment",
    2)
